$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JS-SPA-Self-Evaluation-Protocol")

# "Added edit and delete category" -> mark the related Admin Category rows as
# supported ("Yes") in the score column:
#   Row 43: Admin List Categories
#   Row 44: Admin Create Category
#   Row 47: Admin List Towns
$ws.Range("C43").Value = "Yes"
$ws.Range("C44").Value = "Yes"
$ws.Range("C47").Value = "Yes"

# Update the scroll position / active selection left in the saved view state.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G41").Select()
